$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.616.53'
$ws.Range('E2').Value = '  +1.06%  '
$ws.Range('D3').Value = '3.034.33'
$ws.Range('E3').Value = '  +2.56%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'384.11"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.87%  '
$ws.Range('D6').Value = "'102.50"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.53%  '
$ws.Range('D7').Value = "'0.544"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('D10').Value = "'36.78"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.25%  '
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('D12').Value = "'0.0860"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.80%  '
$ws.Range('D13').Value = '3.516.05'
$ws.Range('E13').Value = '  +2.79%  '
$ws.Range('D14').Value = "'18.58"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.26%  '
$ws.Range('D15').Value = "'7.76"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('D16').Value = '3.048.58'
$ws.Range('E16').Value = '  +2.71%  '
$ws.Range('D17').Value = "'10.84"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -11.83%  '
$ws.Range('E18').Value = '  -3.06%  '
$ws.Range('D19').Value = '51.652.76'
$ws.Range('E19').Value = '  +0.99%  '
$ws.Range('D20').Value = "'3.07"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.74%  '
$ws.Range('D21').Value = "'12.45"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.62%  '
$ws.Range('D22').Value = '0.0₃0962'
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').Value = "'70.03"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.33%  '
$ws.Range('D24').Value = "'267.20"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.46%  '
$ws.Range('D25').Value = "'3.19"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.33%  '
$ws.Range('D26').Value = "'8.34"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.12%  '
$ws.Range('D27').Value = "'7.46"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.94%  '
$ws.Range('E28').Value = '  +3.39%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = "'26.29"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.67%  '
$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').Value = "'1.00"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('E31').Value = '  -1.18%  '
$ws.Range('D32').Value = "'10.29"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.32%  '
$ws.Range('D33').Value = "'2.07"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.20%  '
$ws.Range('D34').Value = "'34.00"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.93%  '
$ws.Range('D35').Value = "'50.50"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.12%  '
$ws.Range('D36').Value = "'0.0446"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.46%  '
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('D39').Value = "'0.291"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.18%  '
$ws.Range('D40').Value = "'17.04"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.38%  '
$ws.Range('E41').Value = '  +1.43%  '
$ws.Range('D42').Value = "'0.116"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('E43').Value = '  +1.56%  '
$ws.Range('D44').Value = "'123.47"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('E45').Value = '  +4.43%  '
$ws.Range('D46').Value = "'21.76"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.68%  '
$ws.Range('E47').Value = '  +2.52%  '
$ws.Range('E48').Value = '  +0.72%  '
$ws.Range('D49').Value = '2.029.31'
$ws.Range('E49').Value = '  -1.43%  '
$ws.Range('D50').Value = '3.337.39'
$ws.Range('E50').Value = '  +2.65%  '
$ws.Range('D51').Value = "'0.518"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.51%  '
